# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45189 (2023-09-20) to 45190 (2023-09-21).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

$range = $ws.Range("C2:C$lastRow")
$range.Value = 45190
